$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2151898734177215
$ws.Range("C2").Value = 0.5158227848101266
$ws.Range("J2").Value = 0.02531645569620253
$ws.Range("P2").Value = 0.1424050632911392
$ws.Range("S2").Value = 0.1012658227848101
$ws.Range("B3").Value = 0.005988023952095809
$ws.Range("C3").Value = 0.03592814371257485
$ws.Range("J3").Value = 0.02994011976047904
$ws.Range("P3").Value = 0.6946107784431138
$ws.Range("S3").Value = 0.2335329341317365
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.673469387755102
$ws.Range("S4").Value = 0.3061224489795918
$ws.Range("B6").Value = 0.03043478260869565
$ws.Range("D6").Value = 0.008695652173913044
$ws.Range("F6").Value = 0.06086956521739131
$ws.Range("J6").Value = 0.2695652173913043
$ws.Range("O6").Value = 0.0391304347826087
$ws.Range("Q6").Value = 0.1434782608695652
$ws.Range("R6").Value = 0.07391304347826087
$ws.Range("S6").Value = 0.3739130434782609
$ws.Range("B7").Value = 0.1268292682926829
$ws.Range("D7").Value = 0.01951219512195122
$ws.Range("F7").Value = 0.03414634146341464
$ws.Range("J7").Value = 0.1463414634146341
$ws.Range("O7").Value = 0.00975609756097561
$ws.Range("Q7").Value = 0.2048780487804878
$ws.Range("R7").Value = 0.08780487804878048
$ws.Range("S7").Value = 0.3707317073170732
$ws.Range("B8").Value = 0.08196721311475409
$ws.Range("D8").Value = 0.02810304449648712
$ws.Range("F8").Value = 0.0585480093676815
$ws.Range("J8").Value = 0.1498829039812646
$ws.Range("O8").Value = 0.00936768149882904
$ws.Range("Q8").Value = 0.1686182669789227
$ws.Range("R8").Value = 0.08196721311475409
$ws.Range("S8").Value = 0.4215456674473068
$ws.Range("B9").Value = 0.1138613861386139
$ws.Range("D9").Value = 0.0198019801980198
$ws.Range("F9").Value = 0.0891089108910891
$ws.Range("J9").Value = 0.1089108910891089
$ws.Range("O9").Value = 0.01485148514851485
$ws.Range("Q9").Value = 0.1584158415841584
$ws.Range("R9").Value = 0.07920792079207921
$ws.Range("S9").Value = 0.4158415841584158
$ws.Range("B10").Value = 0.1171935235158057
$ws.Range("D10").Value = 0.02081727062451812
$ws.Range("F10").Value = 0.07710100231303008
$ws.Range("J10").Value = 0.1133384734001542
$ws.Range("O10").Value = 0.02004626060138782
$ws.Range("Q10").Value = 0.1850424055512722
$ws.Range("R10").Value = 0.06939090208172706
$ws.Range("S10").Value = 0.3970701619121049
$ws.Range("G11").Value = 0.151702786377709
$ws.Range("J11").Value = 0.1114551083591331
$ws.Range("K11").Value = 0.21671826625387
$ws.Range("L11").Value = 0.5139318885448917
$ws.Range("S11").Value = 0.006191950464396285
$ws.Range("G12").Value = 0.7183908045977011
$ws.Range("J12").Value = 0.1954022988505747
$ws.Range("K12").Value = 0.005747126436781609
$ws.Range("L12").Value = 0.04022988505747126
$ws.Range("S12").Value = 0.04022988505747126
$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("F15").Value = 0.01801801801801802
$ws.Range("H15").Value = 0.1261261261261261
$ws.Range("I15").Value = 0.07207207207207207
$ws.Range("J15").Value = 0.3423423423423423
$ws.Range("K15").Value = 0.05855855855855856
$ws.Range("M15").Value = 0.03153153153153153
$ws.Range("O15").Value = 0.03603603603603604
$ws.Range("S15").Value = 0.3153153153153153
$ws.Range("F16").Value = 0.02083333333333333
$ws.Range("H16").Value = 0.15625
$ws.Range("I16").Value = 0.0625
$ws.Range("J16").Value = 0.4010416666666667
$ws.Range("K16").Value = 0.140625
$ws.Range("O16").Value = 0.07291666666666667
$ws.Range("S16").Value = 0.1458333333333333
$ws.Range("F17").Value = 0.02650602409638554
$ws.Range("H17").Value = 0.1783132530120482
$ws.Range("I17").Value = 0.07469879518072289
$ws.Range("J17").Value = 0.4289156626506024
$ws.Range("K17").Value = 0.0963855421686747
$ws.Range("M17").Value = 0.02409638554216868
$ws.Range("O17").Value = 0.06987951807228916
$ws.Range("S17").Value = 0.1012048192771084
$ws.Range("F18").Value = 0.01685393258426966
$ws.Range("H18").Value = 0.1797752808988764
$ws.Range("I18").Value = 0.07865168539325842
$ws.Range("J18").Value = 0.3820224719101123
$ws.Range("K18").Value = 0.0898876404494382
$ws.Range("M18").Value = 0.01685393258426966
$ws.Range("O18").Value = 0.1067415730337079
$ws.Range("S18").Value = 0.1292134831460674
$ws.Range("F19").Value = 0.01574212893553223
$ws.Range("H19").Value = 0.1971514242878561
$ws.Range("I19").Value = 0.09595202398800599
$ws.Range("J19").Value = 0.3748125937031484
$ws.Range("K19").Value = 0.1146926536731634
$ws.Range("M19").Value = 0.01874062968515742
$ws.Range("N19").Value = 0.002998500749625187
$ws.Range("O19").Value = 0.06146926536731634
$ws.Range("S19").Value = 0.1184407796101949
